$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all existing data down by one row (adds a blank row 1, pushing the
# former rows 1-4 to rows 2-5) -- i.e. a column-offset / row-offset insert.
$ws.Rows.Item(1).Insert()

# Row-insert shifts cell values/styles automatically, but the worksheet's
# hyperlink anchors are not re-targeted by the engine, so rebuild them at
# their new locations (D3 and D5), preserving the original address/tooltip/
# display text.
$d3Style = $ws.Range("D3").Style
$d5Style = $ws.Range("D5").Style

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:15866666001@gmail.com", "", "", "15866666001@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:15866666003@gmail.com", "", "mailto:15866666003@gmail.com", "15866666003@gmail.com")

# Adding a hyperlink re-applies the built-in "Hyperlink" cell style; restore
# the original formatting so the cells keep their prior look.
$ws.Range("D3").Style = $d3Style
$ws.Range("D5").Style = $d5Style

# Keep the active selection anchored the same distance below the data as
# before the insert (was E4 of 4 rows; now E10).
$ws.Range("E10").Select()
